# EDD-920: add "Targeted Proteomics Time(s)" and "Metabolomics Time(s)" columns
# to the combinatorial creation sample_experiment_description.xlsx template,
# drop the bold/blue-fill header styling, and refresh the saved window state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
$ws.Range("A1").Value = "Line Name"
$ws.Range("B1").Value = "Replicate Count"
$ws.Range("C1").Value = "Line Description"
$ws.Range("D1").Value = "Part ID"
$ws.Range("E1").Value = "Media"
$ws.Range("F1").Value = "Targeted Proteomics Time(s)"
$ws.Range("G1").Value = "Metabolomics Time(s)"

# --- Data row (row 2) ----------------------------------------------------
$ws.Range("A2").Value = "181-aceF"
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = "Description blah blah"
$ws.Range("D2").Value = "JBx_002078"
$ws.Range("E2").Value = "LB"
$ws.Range("F2").Value = "8h, 24h"
$ws.Range("G2").Value = "4h, 6h"

# --- Remove the old bold-white-on-blue header style (now plain/default) --
$ws.Range("A1:G1").Style = "Normal"

# --- Column widths: drop the old column-E width, size the two new cols ---
$ws.Columns.Item(5).ColumnWidth = 8.43
$ws.Columns.Item(6).ColumnWidth = 16.83203125
$ws.Columns.Item(7).ColumnWidth = 12.666666666666666

# --- Selection / window chrome -------------------------------------------
$ws.Range("D2").Select()
$excel.ActiveWindow.WindowState = -4143

$wb.Windows.Item(1).Left = 0
$wb.Windows.Item(1).Top = 0
